$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update Patient ID, Email, and Past Diagnoses & Treatments
$ws.Range("B2").Value = "H005"
$ws.Range("G2").Value = "test1@gmail.com"
$ws.Range("I2").Value = "Winter is coming"

# Row 3: update Patient ID, Name, Gender, Email, Blood Type, Past Diagnoses & Treatments
$ws.Range("B3").Value = "H025"
$ws.Range("C3").Value = "User21"
$ws.Range("E3").Value = "Female"
$ws.Range("G3").Value = "user3@example.com"
$ws.Range("H3").Value = "B"
$ws.Range("I3").Value = "Fire and Blood"

# Update the active selection to I9 (per the sheetView <selection> change)
$ws.Range("I9").Select() | Out-Null
